# Added Ag and BeO, moved some files from appInfos
$wb = $excel.ActiveWorkbook

# Checklist pattern (columns E..BI) shared by the two new MATERIALS rows -
# identical to the pattern already used on row 18 (Graphite).
$checklist = @(1,0,0,0,0,0,0,0,0,0,1,1,1,1,1,0,0,0,0,0,1,1,1,1,0,0,0,1,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

# --- CERAMICS: add BeO row (row 7) -----------------------------------
# Written first so the shared-string table allocates "BeO" / "Berillium
# oxide" before "Ag" / "Silver".
$wsCeramics = $wb.Worksheets.Item("CERAMICS")
$wsCeramics.Range("A7").Value = 200004
$wsCeramics.Range("B7").Value = "BeO"
$wsCeramics.Range("C7").Value = "BeO"
$wsCeramics.Range("D7").Value = "Berillium oxide"
$wsCeramics.Activate()
$wsCeramics.Range("D7").Select()

# --- METALS AND ALLOYS: add Ag row (row 8) ----------------------------
$wsMetals = $wb.Worksheets.Item("METALS AND ALLOYS")
$wsMetals.Range("A8").Value = 100005
$wsMetals.Range("B8").Value = "Ag"
$wsMetals.Range("C8").Value = "Silver"
$wsMetals.Activate()
$wsMetals.Range("D8").Select()

# --- MATERIALS: add the two new rows (20 = Ag, 21 = BeO) --------------
$wsMaterials = $wb.Worksheets.Item("MATERIALS")

$wsMaterials.Range("A20").Value = 100005
$wsMaterials.Range("B20").Value = "Ag"
$wsMaterials.Range("C20").Value = "Silver"
for ($i = 0; $i -lt $checklist.Length; $i++) {
    $wsMaterials.Cells.Item(20, 5 + $i).Value = $checklist[$i]
}

$wsMaterials.Range("A21").Value = 200005
$wsMaterials.Range("B21").Value = "BeO"
$wsMaterials.Range("C21").Value = "BeO"
for ($i = 0; $i -lt $checklist.Length; $i++) {
    $wsMaterials.Cells.Item(21, 5 + $i).Value = $checklist[$i]
}

$wsMaterials.Activate()
$wsMaterials.Range("D21").Select()
